# Apply updated crypto price/volume data to the worksheet
# (values are textual in the source data, e.g. '26.194.89' or '1.49',
#  so force text number-format before assigning to avoid numeric coercion)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.194.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.592.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.33'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.246'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.02'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.814.79'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.576.13'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.69'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.180.44'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.35'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '213.93'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.04'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.97'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.97'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.423.96'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +8.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.96'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.89%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.588'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.824'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -8.56%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.726.72'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.97'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0956'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.997'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.21%  '
